$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 116 (shifts ipz999 row down to 117)
$ws.Rows.Item(116).Insert()

$ws.Range("A116").Value = "ipz993"
$ws.Range("B116").Value = 99993
$ws.Range("B116").NumberFormat = $ws.Range("B117").NumberFormat
